$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the "Decision Tree" results
$ws.Range("A2").Value = "Decision Tree"
$ws.Range("B2").Value = 0.6294820717131474
$ws.Range("C2").Value = 0.5230769230769231
$ws.Range("D2").Value = 0.7727272727272727
$ws.Range("E2").Value = 0.3953488372093023
$ws.Range("F2").Value = 0.6651416952598805

# Remove the now-obsolete rows (previously Random Forest, K-Nearest Neighbors,
# SVM, Decision Tree, Naive Bayes, XGBoost) so only header + one data row remain
$ws.Range("A3:F8").EntireRow.Delete()

$wb.Save()
